# Add a new "2022" column (column U) to the worksheet, mirroring the
# existing "2021" column (column T) for layout/formatting, and fill in the
# new year's data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the formatting (styles) of column T (rows 3-40, the bordered /
#    formatted block) into the new column U so the new cells inherit the
#    same look (alignment, borders, number format, etc.) as column T.
$ws.Range("T3:T40").Copy()
$ws.Range("U3:U40").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Header row: year label 2022
$ws.Range("U4").Value = 2022

# 3) Data values for the new 2022 column, row by row.
$ws.Range("U6").Value = 1456
$ws.Range("U7").Value = $null

$ws.Range("U8").Value = 45
$ws.Range("U9").Value = 35
$ws.Range("U10").Value = "-"
$ws.Range("U11").Value = 217
$ws.Range("U12").Value = 22
$ws.Range("U13").Value = 8
$ws.Range("U14").Value = "-"
$ws.Range("U15").Value = "-"
$ws.Range("U16").Value = 57
$ws.Range("U17").Value = "-"
$ws.Range("U18").Value = 5
$ws.Range("U19").Value = "-"
$ws.Range("U20").Value = 46
$ws.Range("U21").Value = 1021
$ws.Range("U22").Value = "-"

$ws.Range("U24").Value = 1019
$ws.Range("U25").Value = $null

$ws.Range("U26").Value = 15
$ws.Range("U27").Value = 30
$ws.Range("U28").Value = 1
$ws.Range("U29").Value = 179
$ws.Range("U30").Value = 16
$ws.Range("U31").Value = 8
$ws.Range("U32").Value = "-"
$ws.Range("U33").Value = "-"
$ws.Range("U34").Value = 46
$ws.Range("U35").Value = "-"
$ws.Range("U36").Value = "-"
$ws.Range("U37").Value = "-"
$ws.Range("U38").Value = 25
$ws.Range("U39").Value = 699
$ws.Range("U40").Value = "-"

# 4) Update the view: scroll so column B is the left-most visible column
#    and select V6 (mirrors the author re-saving the file while scrolled).
$ws.Range("B1").Select()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("V6").Select()
